# Terra Nova Template: minor wording update.
# "...is not calculated for Autumn tests." -> "...is not calculated for the Autumn test."
# split into 3 runs (matching the authored edit), and relocate the automatic
# "_GoBack" bookmark onto the second occurrence's new split point.

$d = $word.ActiveDocument

function Split-Run([int]$pos) {
    # Forces a run boundary at $pos (a collapsed Range) without altering any
    # visible content: adding + immediately removing a temp bookmark there
    # splits the underlying run but leaves the text/formatting untouched.
    $tmpName = "ZZSplitTmp"
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($tmpName, $r) | Out-Null
    $d.Bookmarks($tmpName).Delete()
}

function Update-Sentence([bool]$placeBookmark) {
    $searchText = "The National Percentile Rank is not calculated for Autumn tests."

    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }

    $matchStart = $r.Start
    $full = $r.Text

    # Offsets (relative to $matchStart) of the pieces we need to touch.
    $forEndOffset      = $full.IndexOf("for") + 3        # just after "for"
    $autumnStartOffset = $full.IndexOf(" Autumn") + 1    # the "A" of "Autumn"
    $testsStartOffset  = $full.IndexOf("tests")
    $dotOffset         = $full.IndexOf(".", $testsStartOffset)

    # 1) Drop the trailing "s" of "tests" (-> "test").
    $sPos = $matchStart + $dotOffset - 1
    $sRange = $d.Range($sPos, $sPos + 1)
    $sRange.Text = ""

    # 2) Insert "the " right before "Autumn".
    $autumnPos = $matchStart + $autumnStartOffset
    $insPt = $d.Range($autumnPos, $autumnPos)
    $insPt.InsertBefore("the ")

    # New sentence is 3 characters longer than the original ("the " = +4, "s" = -1).
    $newEnd  = $matchStart + $full.Length + 4 - 1
    $split1  = $matchStart + $forEndOffset   # boundary: "...for" | " the Autumn test"
    $split2  = $newEnd - 1                   # boundary: "...test" | "."

    # Re-establish run boundaries. Doing the rightmost boundary first avoids
    # the trailing "." run incorrectly inheriting xml:space="preserve" from
    # its former neighbour.
    Split-Run $split2
    Split-Run $split1
    Split-Run $matchStart

    if ($placeBookmark) {
        # Relocate the document's automatic "_GoBack" bookmark here, between
        # the " the Autumn test" run and the "." run (Bookmarks.Add moves an
        # existing same-named bookmark rather than erroring).
        $bmRange = $d.Range($split2, $split2)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
    }

    return $true
}

Update-Sentence $false | Out-Null
Update-Sentence $true | Out-Null
